$wb = $excel.ActiveWorkbook

# Fix the typo in the first sheet's name
$ws = $wb.Worksheets.Item("Course Equivelents")
$ws.Name = "Course Equivalents"

# Make the (renamed) first sheet the active/selected tab instead of "Grade Schema"
$ws.Activate()
$ws.Select()
